$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: INGRID MENDOZA -> MARLEYDYS CARMONA
$ws.Range("C16").Value = "32907608"
$ws.Range("D16").Value = "MARLEYDYS GREIS CARMONA PADILLA"
$ws.Range("E16").Value = "1908"
$ws.Range("F16").Value = 33125

# Row 17: PAOLA VELEZ -> CESAR MIRANDA
$ws.Range("C17").Value = "1047419950"
$ws.Range("D17").Value = "CESAR LUIS MIRANDA HOYOS"
$ws.Range("E17").Value = "1908"
$ws.Range("G17").Value = 877803

# Row 18: MARLEYDYS CARMONA -> INGRID MENDOZA
$ws.Range("C18").Value = "32939409"
$ws.Range("D18").Value = "INGRID MARGARITA MENDOZA SALAS"
$ws.Range("F18").Value = 37006

# Row 19: MARLEYDYS CARMONA period 1908 -> 1909
$ws.Range("E19").Value = "1909"

# Row 20 unchanged (CESAR MIRANDA, 1909, 33125, 877803)

# Row 21: CESAR MIRANDA -> PAOLA VELEZ
$ws.Range("C21").Value = "45554716"
$ws.Range("D21").Value = "PAOLA ROCIO VELEZ ULLOQUE"
$ws.Range("E21").Value = "1911"
$ws.Range("G21").Value = 828116
